$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - copy style from an existing header cell (H1) to keep the same
# bold/centered/bordered formatting used by the other headers.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the two new columns.
$values = @{
    2 = 7
    3 = 9
    4 = 6
    5 = 9
    6 = 8
    7 = 7
    8 = 4
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Range("I$row").Value = $v
    $ws.Range("J$row").Value = $v
}
